$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style used by the
# other header cells (bold, centered, bordered) by copying formatting from H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells I2 and J2 (plain, unstyled, like the other row 2 values)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9
